# "homogenized D1S input library" -- adds a new FNG experimental benchmark
# row and bumps the Oktavian NPS cut-off on the "Experimental benchmarks"
# sheet, and makes that sheet the active one.

$wb = $excel.ActiveWorkbook

$wsExp = $wb.Worksheets.Item("Experimental benchmarks")

# Bump the NPS cut-off for the existing Oktavian experiment row.
$wsExp.Range("F4").Value = 1000000

# New row 5: Frascati Neutron Generator (FNG) SDDR experiment.
$wsExp.Range("A5").Value = "Frascati Neutron Generator SDDR experimen"
$wsExp.Range("B5").Value = "FNG"

# C5/D5/E5 need to hold the literal text "false" (same shared string used
# by C4), not a COM boolean -- copy/paste-values from C4 keeps it typed as
# text instead of auto-coercing the string "false" into TRUE/FALSE.
$wsExp.Range("C4").Copy()
$wsExp.Range("C5").PasteSpecial(-4163)
$wsExp.Range("C4").Copy()
$wsExp.Range("D5").PasteSpecial(-4163)
$wsExp.Range("C4").Copy()
$wsExp.Range("E5").PasteSpecial(-4163)

$wsExp.Range("F5").Value = 500000000
$wsExp.Range("J5").Value = "D1S5"

# Make "Experimental benchmarks" the active sheet/selection (was previously
# "Computational benchmarks"), with F5 selected.
$wsExp.Activate()
$wsExp.Range("F5").Select()
